$d = $word.ActiveDocument

# The footnote with id="2" (the only "real" footnote - Footnotes.Item(1) in the
# collection, since separator/continuationSeparator footnotes aren't exposed
# there) gets:
#   - a new <w:pPr><w:tabs>...</w:tabs></w:pPr> on its paragraph
#   - a new run containing <w:tab/> before the footnoteRef run
#   - a new run containing <w:tab/> after the footnoteRef run
#   - the leading space removed from the text run
#
# Footnote.Range's InsertBefore/InsertAfter only let us touch the very start
# or very end of the footnote paragraph, so instead we replace the whole
# paragraph's contents in one shot via InsertXML, using the
# pkg:package/WordOpenXML envelope so Word treats it as a full paragraph
# replacement rather than trying to merge inline text.

$fn = $d.Footnotes.Item(1)
$r = $fn.Range

$xml = '<?xml version="1.0"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
'<w:p>' +
'<w:pPr><w:tabs><w:tab w:val="start" w:pos="0"/><w:tab w:val="start" w:pos="400"/></w:tabs></w:pPr>' +
'<w:r><w:tab/></w:r>' +
'<w:r><w:rPr><w:shadow/><w:vertAlign w:val="superscript"/></w:rPr><w:footnoteRef/></w:r>' +
'<w:r><w:tab/></w:r>' +
'<w:r><w:t xml:space="preserve">This is the content of the footnote.</w:t></w:r>' +
'</w:p>' +
'</w:body>' +
'</w:document>' +
'</pkg:xmlData>' +
'</pkg:part>' +
'</pkg:package>'

$r.InsertXML($xml)
